$wb = $excel.ActiveWorkbook

# Delete column H ("是否有舞台（字符串匹配）") on every sheet.
# This shifts the old "Link" column (I) into H, and the old "Cover"
# column (J) into I, matching the new layout.
foreach ($ws in $wb.Worksheets) {
    $ws.Columns.Item(8).Delete()
}

# Sheet "展览" (1st sheet): fix the data-type/value regressions that
# happened on top of the column shift.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value = 227
$ws1.Cells.Item(2, 7).Value = 45
$ws1.Cells.Item(3, 7).Value = 50
$ws1.Cells.Item(4, 7).Value = 45
$ws1.Cells.Item(5, 7).Value = 45

# Sheet "全部类型" (4th sheet): same fix-ups as sheet 1.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 227
$ws4.Cells.Item(2, 7).Value = 45
$ws4.Cells.Item(3, 7).Value = 50
$ws4.Cells.Item(4, 7).Value = 45
$ws4.Cells.Item(5, 7).Value = 45
